# API: Gameweeks import (#25)
# Adds two new columns to the "Challenges" sheet: "Show Statistics Continuously"
# and "Gameweek", populated for both existing challenge rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Challenges")

# New header cells
$ws.Range("S1").Value = "Show Statistics Continuously"
$ws.Range("T1").Value = "Gameweek"

# Row 2 (C001) - "true" must be stored as literal text, not a boolean
$ws.Range("S2").Value = "'true"
$ws.Range("S2").Style = "Normal"
$ws.Range("T2").Value = 1

# Row 3 (C002)
$ws.Range("S3").Value = "'true"
$ws.Range("S3").Style = "Normal"
$ws.Range("T3").Value = 2
